$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 640.5909
$ws.Range("I28").Value = 586.6875
$ws.Range("J28").Value = 784.3333
$ws.Range("K28").Value = 586.6875
$ws.Range("L28").Value = 784.3333
$ws.Range("M28").Value = -101.6875
$ws.Range("N28").Value = -1754.3333
$ws.Range("M98").Value = 325.0999999999999
$ws.Range("M122").Value = -1068.7
$ws.Range("H129").Value = 891.8095
$ws.Range("J129").Value = 905.0877
$ws.Range("L129").Value = 2715.2631
$ws.Range("N129").Value = -12715.2631
$ws.Range("H137").Value = 1998
$ws.Range("I137").Value = 1497.5
$ws.Range("K137").Value = 4492.5
$ws.Range("M137").Value = -1942.5
$ws.Range("N98").ClearContents()
$ws.Range("N122").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3723.5557
$ws.Range("I61").Value = 2824.6743
$ws.Range("J61").Value = 7237.364
$ws.Range("K61").Value = 2824.6743
$ws.Range("L61").Value = 7237.364
$ws.Range("M61").Value = -2612.6743
$ws.Range("N61").Value = -7661.364
$ws.Range("H74").Value = 2330.8
$ws.Range("I74").Value = 2580.2
$ws.Range("J74").Value = 2081.4
$ws.Range("K74").Value = 2580.2
$ws.Range("L74").Value = 2081.4
$ws.Range("M74").Value = -1706.2
$ws.Range("N74").Value = -3829.4
$ws.Range("H77").Value = 2330.8
$ws.Range("I77").Value = 2580.2
$ws.Range("J77").Value = 2081.4
$ws.Range("K77").Value = 12901
$ws.Range("L77").Value = 10407
$ws.Range("M77").Value = -8533
$ws.Range("N77").Value = -19143
$ws.Range("H132").Value = 3138.56
$ws.Range("I132").Value = 1533.6923
$ws.Range("J132").Value = 4877.1665
$ws.Range("K132").Value = 4601.0769
$ws.Range("L132").Value = 14631.4995
$ws.Range("M132").Value = -2071.0769
$ws.Range("N132").Value = -19691.4995
$ws.Range("H136").Value = 3723.5557
$ws.Range("I136").Value = 2824.6743
$ws.Range("J136").Value = 7237.364
$ws.Range("K136").Value = 8474.0229
$ws.Range("L136").Value = 21712.092
$ws.Range("M136").Value = -5924.0229
$ws.Range("N136").Value = -26812.092

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3866.3076
$ws.Range("I134").Value = 2424.8
$ws.Range("J134").Value = 8671.333000000001
$ws.Range("K134").Value = 7274.400000000001
$ws.Range("L134").Value = 26013.999
$ws.Range("M134").Value = -4739.400000000001
$ws.Range("N134").Value = -31083.999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15527.4
$ws.Range("I31").Value = 3612
$ws.Range("J31").Value = 25276.363
$ws.Range("K31").Value = 3612
$ws.Range("L31").Value = 25276.363
$ws.Range("M31").Value = -3317
$ws.Range("N31").Value = -25866.363
$ws.Range("H34").Value = 15527.4
$ws.Range("I34").Value = 3612
$ws.Range("J34").Value = 25276.363
$ws.Range("K34").Value = 3612
$ws.Range("L34").Value = 25276.363
$ws.Range("M34").Value = -3410
$ws.Range("N34").Value = -25680.363
$ws.Range("H58").Value = 2014.9445
$ws.Range("I58").Value = 1534.5
$ws.Range("J58").Value = 2255.1667
$ws.Range("K58").Value = 1534.5
$ws.Range("L58").Value = 2255.1667
$ws.Range("M58").Value = -1331.5
$ws.Range("N58").Value = -2661.1667
$ws.Range("H132").Value = 2732.25
$ws.Range("I132").Value = 2128.1
$ws.Range("K132").Value = 6384.299999999999
$ws.Range("M132").Value = -3854.299999999999
$ws.Range("H134").Value = 18004.8
$ws.Range("I134").Value = 22670
$ws.Range("J134").Value = 11007
$ws.Range("K134").Value = 68010
$ws.Range("L134").Value = 33021
$ws.Range("M134").Value = -65475
$ws.Range("N134").Value = -38091
$ws.Range("H136").Value = 2014.9445
$ws.Range("I136").Value = 1534.5
$ws.Range("J136").Value = 2255.1667
$ws.Range("K136").Value = 4603.5
$ws.Range("L136").Value = 6765.500100000001
$ws.Range("M136").Value = -2053.5
$ws.Range("N136").Value = -11865.5001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1090.7667
$ws.Range("I5").Value = 576.1429000000001
$ws.Range("J5").Value = 1541.0625
$ws.Range("K5").Value = 1728.4287
$ws.Range("L5").Value = 4623.1875
$ws.Range("M5").Value = -1616.4287
$ws.Range("N5").Value = -4847.1875
$ws.Range("H131").Value = 1110
$ws.Range("J131").Value = 1240.8823
$ws.Range("L131").Value = 3722.6469
$ws.Range("N131").Value = -13802.6469
$ws.Range("H135").Value = 1090.7667
$ws.Range("I135").Value = 576.1429000000001
$ws.Range("J135").Value = 1541.0625
$ws.Range("K135").Value = 5185.2861
$ws.Range("L135").Value = 13869.5625
$ws.Range("M135").Value = -2650.2861
$ws.Range("N135").Value = -18939.5625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1765.975
$ws.Range("I126").Value = 1456.8148
$ws.Range("J126").Value = 2408.077
$ws.Range("K126").Value = 4370.4444
$ws.Range("L126").Value = 7224.231000000001
$ws.Range("M126").Value = -1900.4444
$ws.Range("N126").Value = -12164.231

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3527.0715
$ws.Range("I122").Value = 2361.625
$ws.Range("J122").Value = 5081
$ws.Range("K122").Value = 7084.875
$ws.Range("L122").Value = 15243
$ws.Range("M122").Value = -4634.875
$ws.Range("N122").Value = -20143
$ws.Range("H132").Value = 6918.3394
$ws.Range("I132").Value = 9184.1875
$ws.Range("J132").Value = 3897.2083
$ws.Range("K132").Value = 27552.5625
$ws.Range("L132").Value = 11691.6249
$ws.Range("M132").Value = -25022.5625
$ws.Range("N132").Value = -16751.6249
$ws.Range("H136").Value = 5682.357
$ws.Range("I136").Value = 4754.8
$ws.Range("J136").Value = 8001.25
$ws.Range("K136").Value = 14264.4
$ws.Range("L136").Value = 24003.75
$ws.Range("M136").Value = -11714.4
$ws.Range("N136").Value = -29103.75
$ws.Range("H140").Value = 67483.39999999999
$ws.Range("J140").Value = 67483.39999999999
$ws.Range("L140").Value = 67483.39999999999
$ws.Range("N140").Value = -77843.39999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2124.5386
$ws.Range("I132").Value = 1509.5769
$ws.Range("J132").Value = 3354.4614
$ws.Range("K132").Value = 4528.7307
$ws.Range("L132").Value = 10063.3842
$ws.Range("M132").Value = -1998.7307
$ws.Range("N132").Value = -15123.3842
$ws.Range("H136").Value = 8158782.5
$ws.Range("I136").Value = 25718512
$ws.Range("J136").Value = 6050.7144
$ws.Range("K136").Value = 77155536
$ws.Range("L136").Value = 18152.1432
$ws.Range("M136").Value = -77152986
$ws.Range("N136").Value = -23252.1432
